# Replace each two-digit multiplication problem text in the table
# with its new randomized equivalent, per the source diff.
$d = $word.ActiveDocument

$d.Content.Find.Execute("51×47=", $true, $false, $false, $false, $false, $true, 1, $false, "23×42=", 2) | Out-Null
$d.Content.Find.Execute("59×86=", $true, $false, $false, $false, $false, $true, 1, $false, "44×51=", 2) | Out-Null
$d.Content.Find.Execute("53×36=", $true, $false, $false, $false, $false, $true, 1, $false, "98×55=", 2) | Out-Null
$d.Content.Find.Execute("16×70=", $true, $false, $false, $false, $false, $true, 1, $false, "45×50=", 2) | Out-Null
$d.Content.Find.Execute("91×65=", $true, $false, $false, $false, $false, $true, 1, $false, "47×78=", 2) | Out-Null
$d.Content.Find.Execute("69×61=", $true, $false, $false, $false, $false, $true, 1, $false, "17×60=", 2) | Out-Null
$d.Content.Find.Execute("98×23=", $true, $false, $false, $false, $false, $true, 1, $false, "53×35=", 2) | Out-Null
$d.Content.Find.Execute("66×63=", $true, $false, $false, $false, $false, $true, 1, $false, "33×77=", 2) | Out-Null
$d.Content.Find.Execute("39×44=", $true, $false, $false, $false, $false, $true, 1, $false, "34×41=", 2) | Out-Null
$d.Content.Find.Execute("94×72=", $true, $false, $false, $false, $false, $true, 1, $false, "98×27=", 2) | Out-Null
$d.Content.Find.Execute("51×11=", $true, $false, $false, $false, $false, $true, 1, $false, "98×26=", 2) | Out-Null
$d.Content.Find.Execute("27×64=", $true, $false, $false, $false, $false, $true, 1, $false, "64×29=", 2) | Out-Null
$d.Content.Find.Execute("59×99=", $true, $false, $false, $false, $false, $true, 1, $false, "87×52=", 2) | Out-Null
$d.Content.Find.Execute("23×36=", $true, $false, $false, $false, $false, $true, 1, $false, "59×58=", 2) | Out-Null
$d.Content.Find.Execute("52×88=", $true, $false, $false, $false, $false, $true, 1, $false, "87×81=", 2) | Out-Null
$d.Content.Find.Execute("56×48=", $true, $false, $false, $false, $false, $true, 1, $false, "83×56=", 2) | Out-Null
$d.Content.Find.Execute("53×47=", $true, $false, $false, $false, $false, $true, 1, $false, "20×82=", 2) | Out-Null
$d.Content.Find.Execute("24×92=", $true, $false, $false, $false, $false, $true, 1, $false, "64×55=", 2) | Out-Null
$d.Content.Find.Execute("19×31=", $true, $false, $false, $false, $false, $true, 1, $false, "99×30=", 2) | Out-Null
$d.Content.Find.Execute("86×83=", $true, $false, $false, $false, $false, $true, 1, $false, "60×43=", 2) | Out-Null
$d.Content.Find.Execute("14×30=", $true, $false, $false, $false, $false, $true, 1, $false, "33×47=", 2) | Out-Null
$d.Content.Find.Execute("75×66=", $true, $false, $false, $false, $false, $true, 1, $false, "25×93=", 2) | Out-Null
$d.Content.Find.Execute("40×28=", $true, $false, $false, $false, $false, $true, 1, $false, "47×34=", 2) | Out-Null
$d.Content.Find.Execute("69×16=", $true, $false, $false, $false, $false, $true, 1, $false, "19×75=", 2) | Out-Null
$d.Content.Find.Execute("71×74=", $true, $false, $false, $false, $false, $true, 1, $false, "87×51=", 2) | Out-Null
